$d = $word.ActiveDocument

# The document currently ends with the "Dr Emma Kilford <...>" paragraph,
# immediately followed by the section properties. We need to append, after
# that paragraph:
#   1) a new paragraph for "Dr Laura Panagi <lp579@medschl.cam.ac.uk>"
#   2) three new, empty paragraphs
# all of which carry the same "iCs" (italic - complex scripts) paragraph
# mark formatting used throughout the supervisors list.

$lastPara = $d.Paragraphs.Last
$tail = $lastPara.Range

# Create a fresh paragraph right after the last one; its Range (the new
# paragraph mark) is what we will overwrite with the real OOXML content so
# that no stray empty runs are left behind.
$tail.InsertParagraphAfter() | Out-Null
$newRange = $d.Paragraphs.Last.Range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:pPr><w:rPr><w:iCs/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:iCs/></w:rPr><w:t xml:space="preserve">Dr Laura </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:iCs/></w:rPr><w:t>Panagi</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:iCs/></w:rPr><w:t xml:space="preserve"> &lt;lp579@medschl.cam.ac.uk&gt;</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:rPr><w:iCs/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:rPr><w:iCs/></w:rPr></w:pPr></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($xml)
